$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44476
$ws.Range("I2").Value = '1a nueva(o)'
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 490
$ws.Range("O2").Value = 'Perú'
$ws.Range("P2").Value = 490
$ws.Range("D4").Value = 44349
$ws.Range("H4").Value = 'Pachia'
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 730
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = 740
$ws.Range("P4").Value = 740
$ws.Range("D5").Value = 44670
$ws.Range("K5").Value = 400
$ws.Range("L5").Value = 430
$ws.Range("M5").Value = 415
$ws.Range("P5").Value = 415
$ws.Range("D6").Value = 44620
$ws.Range("I6").Value = '1a (cosecha)'
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 480
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = 490
$ws.Range("P6").Value = 490
$ws.Range("D7").Value = 44211
$ws.Range("J7").Value = 1600
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 550
$ws.Range("M7").Value = 525
$ws.Range("O7").Value = 'Región de O''Higgins'
$ws.Range("P7").Value = 525
$ws.Range("D8").Value = 44490
$ws.Range("I8").Value = '1a nueva(o)'
$ws.Range("K8").Value = 450
$ws.Range("L8").Value = 480
$ws.Range("M8").Value = 465
$ws.Range("O8").Value = 'Perú'
$ws.Range("P8").Value = 465
$ws.Range("D9").Value = 44630
$ws.Range("K9").Value = 450
$ws.Range("L9").Value = 480
$ws.Range("M9").Value = 465
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 465
$ws.Range("D10").Value = 44243
$ws.Range("I10").Value = '2a nueva(o)'
$ws.Range("J10").Value = 1600
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("D11").Value = 44575
$ws.Range("I11").Value = '1a nueva(o)'
$ws.Range("J11").Value = 1300
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 550
$ws.Range("M11").Value = 525
$ws.Range("O11").Value = 'Región de O''Higgins'
$ws.Range("P11").Value = 525
$ws.Range("D12").Value = 44322
$ws.Range("I12").Value = '1a (cosecha)'
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = 375
$ws.Range("O12").Value = 'Región del Maule'
$ws.Range("P12").Value = 375
$ws.Range("D13").Value = 44530
$ws.Range("J13").Value = 900
$ws.Range("K13").Value = 480
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = 490
$ws.Range("O13").Value = 'Perú'
$ws.Range("P13").Value = 490
$ws.Range("D14").Value = 44201
$ws.Range("I14").Value = '1a nueva(o)'
$ws.Range("J14").Value = 1360
$ws.Range("K14").Value = 730
$ws.Range("L14").Value = 750
$ws.Range("M14").Value = 740
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 740
$ws.Range("D15").Value = 44665
$ws.Range("I15").Value = '1a (cosecha)'
$ws.Range("J15").Value = 1200
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 420
$ws.Range("M15").Value = 410
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 410
$ws.Range("D16").Value = 44179
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1350
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = 1375
$ws.Range("P16").Value = 1375
$ws.Range("D17").Value = 44428
$ws.Range("I17").Value = '1a nueva(o)'
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 580
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = 590
$ws.Range("O17").Value = 'Perú'
$ws.Range("P17").Value = 590
$ws.Range("D18").Value = 44691
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 580
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = 590
$ws.Range("P18").Value = 590
$ws.Range("D19").Value = 44601
$ws.Range("I19").Value = '2a (cosecha)'
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 450
$ws.Range("M19").Value = 425
$ws.Range("P19").Value = 425
$ws.Range("D20").Value = 44736
$ws.Range("I20").Value = '1a (guarda)'
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 780
$ws.Range("L20").Value = 800
$ws.Range("M20").Value = 790
$ws.Range("O20").Value = 'Región de O''Higgins'
$ws.Range("P20").Value = 790
$ws.Range("D21").Value = 44685
$ws.Range("I21").Value = '1a (cosecha)'
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 680
$ws.Range("L21").Value = 700
$ws.Range("M21").Value = 690
$ws.Range("O21").Value = 'Región de O''Higgins'
$ws.Range("P21").Value = 690
$ws.Range("D22").Value = 44483
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 580
$ws.Range("M22").Value = 565
$ws.Range("O22").Value = 'Perú'
$ws.Range("P22").Value = 565
$ws.Range("D23").Value = 44175
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1550
$ws.Range("M23").Value = 1525
$ws.Range("P23").Value = 1525
$ws.Range("D24").Value = 44469
$ws.Range("J24").Value = 1200
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 650
$ws.Range("M24").Value = 625
$ws.Range("O24").Value = 'Perú'
$ws.Range("P24").Value = 625
$ws.Range("D25").Value = 44301
$ws.Range("H25").Value = 'Camote'
$ws.Range("I25").Value = '2a nueva(o)'
$ws.Range("K25").Value = 400
$ws.Range("L25").Value = 430
$ws.Range("M25").Value = 415
$ws.Range("O25").Value = 'Provincia de Melipilla'
$ws.Range("P25").Value = 415
$ws.Range("D26").Value = 44204
$ws.Range("I26").Value = '2a nueva(o)'
$ws.Range("J26").Value = 1600
$ws.Range("K26").Value = 500
$ws.Range("L26").Value = 550
$ws.Range("M26").Value = 525
$ws.Range("P26").Value = 525
$ws.Range("D27").Value = 44231
$ws.Range("I27").Value = '1a nueva(o)'
$ws.Range("J27").Value = 1300
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 480
$ws.Range("M27").Value = 465
$ws.Range("P27").Value = 465
$ws.Range("D28").Value = 44238
$ws.Range("J28").Value = 1250
$ws.Range("K28").Value = 430
$ws.Range("L28").Value = 450
$ws.Range("M28").Value = 440
$ws.Range("P28").Value = 440
$ws.Range("D29").Value = 44399
$ws.Range("I29").Value = '1a (guarda)'
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 450
$ws.Range("L29").Value = 480
$ws.Range("M29").Value = 465
$ws.Range("O29").Value = 'Provincia de Melipilla'
$ws.Range("P29").Value = 465
$ws.Range("D30").Value = 44547
$ws.Range("I30").Value = '1a (cosecha)'
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 650
$ws.Range("M30").Value = 625
$ws.Range("P30").Value = 625
$ws.Range("I31").Value = '2a nueva(o)'
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 550
$ws.Range("M31").Value = 525
$ws.Range("P31").Value = 525
$ws.Range("D32").Value = 44453
$ws.Range("I32").Value = '1a nueva(o)'
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 630
$ws.Range("L32").Value = 650
$ws.Range("M32").Value = 640
$ws.Range("P32").Value = 640
